$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "5XDOPW"
$ws.Range("B10").Value = "Almohadilla + Chip Epson C9345"
$ws.Range("C10").Value = "L6550 L6570 L6580 L8050 L8160 L8168 L8180 L8188 L11160 L15140 L15150 L15158, L15160 L15168 L15180 L18150, ET 5800 5150 5850 5880 16150 16600 16650, WF 3820 4820 4830 7800 7820 7830 7840 7845, ET M15140 M16600, ST C8000, EC C7000,  ITS L6550 L6570 L6580, PX M6010F M6011F M6711FT M6712FT  M791FT S6710T, EW M873, TEW M973A3"
$ws.Range("D10").Value = 35000
$ws.Range("E10").Value = 200000
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 8
$ws.Range("H10").Formula = "=(E10-D10)*G10"
$ws.Range("I10").Formula = "=D10*F10"
$ws.Range("J10").Value = 0
